$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Signature proportions")
$ws.Columns.Item(1).Delete()
